# Update automàtic: dades i banners [2026-02-28 20:19]
# Refresh DATA_EXTRACCIO timestamps and the corresponding observation values
# pulled from meteo.cat for the affected stations/rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-28 20:18:17"
$ws.Range("I2").Value = "0.2 mm"
$ws.Range("O2").Value = "2.9 °C"
$ws.Range("E3").Value = "2026-02-28 20:18:20"
$ws.Range("I3").Value = "0.4 mm"
$ws.Range("N3").Value = "-3.0 °C 19:51 TU"
$ws.Range("E4").Value = "2026-02-28 20:18:22"
$ws.Range("J4").Value = "1024.8 hPa"
$ws.Range("O4").Value = "11.1 °C"
$ws.Range("E5").Value = "2026-02-28 20:18:25"
$ws.Range("N5").Value = "-2.9 °C 19:55 TU"
$ws.Range("O5").Value = "-1.3 °C"
$ws.Range("E6").Value = "2026-02-28 20:18:27"
$ws.Range("J6").Value = "1024.7 hPa"
$ws.Range("N6").Value = "9.8 °C 19:58 TU"
$ws.Range("O6").Value = "12.0 °C"
$ws.Range("E7").Value = "2026-02-28 20:18:30"
$ws.Range("J7").Value = "1024.8 hPa"
$ws.Range("E8").Value = "2026-02-28 20:18:32"
$ws.Range("E9").Value = "2026-02-28 20:18:35"
$ws.Range("H9").Value = "'79%"
$ws.Range("E10").Value = "2026-02-28 20:18:36"
$ws.Range("O10").Value = "10.9 °C"
$ws.Range("E11").Value = "2026-02-28 20:18:37"
$ws.Range("E12").Value = "2026-02-28 20:18:38"
$ws.Range("H12").Value = "'82%"
$ws.Range("E13").Value = "2026-02-28 20:18:39"
$ws.Range("E14").Value = "2026-02-28 20:18:40"
$ws.Range("H14").Value = "'83%"
$ws.Range("E15").Value = "2026-02-28 20:18:41"
$ws.Range("O15").Value = "11.1 °C"
$ws.Range("E16").Value = "2026-02-28 20:18:42"
$ws.Range("H16").Value = "'66%"
$ws.Range("E17").Value = "2026-02-28 20:18:43"
$ws.Range("H17").Value = "'82%"
$ws.Range("N17").Value = "1.4 °C 19:50 TU"
$ws.Range("E18").Value = "2026-02-28 20:18:44"
$ws.Range("O18").Value = "11.7 °C"
$ws.Range("E19").Value = "2026-02-28 20:18:45"
$ws.Range("E20").Value = "2026-02-28 20:18:47"
$ws.Range("H20").Value = "'63%"
$ws.Range("N20").Value = "-2.0 °C 19:59 TU"
$ws.Range("E21").Value = "2026-02-28 20:18:49"
$ws.Range("E22").Value = "2026-02-28 20:18:51"
$ws.Range("E23").Value = "2026-02-28 20:18:54"
$ws.Range("I23").Value = "0.3 mm"
$ws.Range("N23").Value = "-2.3 °C 19:54 TU"
$ws.Range("E24").Value = "2026-02-28 20:18:56"
$ws.Range("E25").Value = "2026-02-28 20:18:58"
$ws.Range("H25").Value = "'62%"
$ws.Range("N25").Value = "-1.1 °C 19:54 TU"
$ws.Range("E26").Value = "2026-02-28 20:19:01"
$ws.Range("E27").Value = "2026-02-28 20:19:03"
$ws.Range("H27").Value = "'55%"
$ws.Range("N27").Value = "-0.5 °C 19:49 TU"
$ws.Range("O27").Value = "1.9 °C"
$ws.Range("E28").Value = "2026-02-28 20:19:06"
$ws.Range("E29").Value = "2026-02-28 20:19:08"
$ws.Range("H29").Value = "'84%"
$ws.Range("K29").Value = "12.3 MJ/m2"
$ws.Range("E30").Value = "2026-02-28 20:19:11"
$ws.Range("E31").Value = "2026-02-28 20:19:13"
$ws.Range("H31").Value = "'81%"
$ws.Range("L31").Value = "66.2 km/h - 339º 19:57 TU"
$ws.Range("E32").Value = "2026-02-28 20:19:16"
$ws.Range("E33").Value = "2026-02-28 20:19:18"
$ws.Range("H33").Value = "'65%"
$ws.Range("J33").Value = "1023.0 hPa"
$ws.Range("O33").Value = "7.2 °C"
$ws.Range("E34").Value = "2026-02-28 20:19:21"
$ws.Range("H34").Value = "'68%"
$ws.Range("I34").Value = "0.4 mm"
$ws.Range("E35").Value = "2026-02-28 20:19:23"
$ws.Range("E36").Value = "2026-02-28 20:19:26"
$ws.Range("E37").Value = "2026-02-28 20:19:28"
$ws.Range("J37").Value = "1025.9 hPa"
$ws.Range("E38").Value = "2026-02-28 20:19:31"
$ws.Range("O38").Value = "11.7 °C"
$ws.Range("E39").Value = "2026-02-28 20:19:33"
$ws.Range("H39").Value = "'63%"
$ws.Range("N39").Value = "-2.0 °C 19:33 TU"
$ws.Range("E40").Value = "2026-02-28 20:19:35"
$ws.Range("E41").Value = "2026-02-28 20:19:38"
$ws.Range("H41").Value = "'72%"
$ws.Range("L41").Value = "14.8 km/h - 29º 19:55 TU"
$ws.Range("E42").Value = "2026-02-28 20:19:40"
$ws.Range("E43").Value = "2026-02-28 20:19:42"
$ws.Range("O43").Value = "7.7 °C"
$ws.Range("E44").Value = "2026-02-28 20:19:45"
$ws.Range("H44").Value = "'91%"
$ws.Range("O44").Value = "-1.1 °C"
$ws.Range("E45").Value = "2026-02-28 20:19:47"
$ws.Range("N45").Value = "4.3 °C 19:56 TU"
$ws.Range("E46").Value = "2026-02-28 20:19:50"
$ws.Range("N46").Value = "9.1 °C 19:55 TU"
$ws.Range("O46").Value = "11.5 °C"
